# Workbook was edited to add a new daily-routine log entry (row 10) and
# update the active selection/scroll position on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row 10 entry -------------------------------------------------
# Copy the date formatting from the previous date cell (A9) so the new date
# cell (A10) reuses the existing "short date" style instead of creating a
# brand-new number-format style.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A10").Value = 43816         # 12/17/2019

$ws.Range("B10").Value = "Edit the document"
$ws.Range("C10").Value = "Feasibility and DFD"

# --- Update the view: scroll/selection -----------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3   # topLeftCell column -> "C"
$win.ScrollRow = 1      # topLeftCell row -> 1
$ws.Range("D12").Select()
